# Small poster edit:
#  - Merge the two-line title into a single line.
#  - Reposition/resize the Abstract box and append a clause about the
#    stationary-accelerometer reset prompt to its last bullet.
#  - Nudge a few neighbouring boxes/pictures up and to the side to close
#    the gap left by the now-shorter title / taller abstract box.

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title ("Hybrid Localization with Video Based Positioning" / "Technology") ---
# Collapse the two paragraphs into one; the textbox has spAutoFit so its
# height shrinks to match automatically.
$titleShape = Get-ShapeById $s 32
$titleShape.TextFrame.TextRange.Text = "Hybrid Localization with Video Based Positioning Technology"

# --- Abstract box: move up, grow taller, extend the last bullet's text ---
$abstractShape = Get-ShapeById $s 33
$abstractText = $abstractShape.TextFrame.TextRange
$lastBullet = $abstractText.Characters(395, 254)
$lastBullet.Text = "This proposed solution provides for an infrastructure-less and low cost solution in the indoor localization field by applying pedestrian dead reckoning using the smartphone’s sensors and camera to track the user’s location as well as to scan the QR Code as well as a stationary accelerometer prompt to allow users to reset the accelerometer."
$abstractShape.Top = 225.5502471923828
$abstractShape.Height = 490.4360046386719

# --- "Technologies/Tools Used" box: follow the abstract box upward ---
$techShape = Get-ShapeById $s 15
$techShape.Top = 225.5502471923828

# --- "System Architecture" heading: shift down slightly ---
$sysArchShape = Get-ShapeById $s 18
$sysArchShape.Top = 686.1734008789062

# --- "Figure 1. Indoor Localization application architecture" caption ---
$fig1Shape = Get-ShapeById $s 10
$fig1Shape.Top = 1044.705322265625

# --- Sensor description body text box ---
$descrShape = Get-ShapeById $s 24
$descrShape.Left = 760.877685546875
$descrShape.Top = 641.8773803710938

# --- Picture 70 (the "System Architecture" diagram image) ---
$picShape = Get-ShapeById $s 71
$picShape.Top = 758.31103515625
